# Updated cryptos list on Sun Feb  4 05:37:07 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like plain numbers (e.g. "97.46").
# Force those specific cells to Text format before writing so Excel keeps them
# as literal strings instead of silently parsing them into numeric values.

$ws.Range("D2").Value = "42.903.68"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.302.72"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.81"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.46"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.77"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.75"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").Value = "2.659.85"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "2.301.08"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "42.885.74"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.59"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.94"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.13"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.00"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.22"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.02"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.63"
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.42"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0685"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "2.000.31"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.28"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "2.524.47"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.24"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.18"
$ws.Range("E51").Value = "  -5.34%  "
